$d = $word.ActiveDocument

# Merge split runs (with stray proofErr spell-check markers) for the
# professor's name back into a single contiguous run.
$d.Content.Find.Execute(
    "Professor(a):  Miriã da Silveira Coelho Corrêa     ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Professor(a):  Miriã da Silveira Coelho Corrêa     ", 2)

# Merge split runs (with stray proofErr grammar-check markers) for the
# late-submission sentence back into a single contiguous run.
$d.Content.Find.Execute(
    "Caso o trabalho seja enviado após a data estipulada, o mesmo será pontuado com 0.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Caso o trabalho seja enviado após a data estipulada, o mesmo será pontuado com 0.", 2)

# Merge the split "menu com Links" run back into a single run.
$d.Content.Find.Execute(
    "Imagem mapeada contendo o menu com Links",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Imagem mapeada contendo o menu com Links", 2)
